$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 266; this shifts rows 266-377 down to 267-378
$ws.Rows(266).Insert()

# Populate the newly inserted row 266 with the new record's data
$ws.Range("A266").Value = 9
$ws.Range("B266").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C266").Value = "Metropolitana"
$ws.Range("D266").Value = 44704
$ws.Range("E266").Value = 13
$ws.Range("F266").Value = 100112044
$ws.Range("G266").Value = "Perejil"
$ws.Range("H266").Value = "Sin especificar"
$ws.Range("I266").Value = "Primera"
$ws.Range("J266").Value = 25
$ws.Range("K266").Value = 10000
$ws.Range("L266").Value = 12000
$ws.Range("M266").Value = 11200
$ws.Range("N266").Value = "$/docena de atados"
$ws.Range("O266").Value = "Región Metropolitana"
$ws.Range("P266").Value = 3733
$ws.Range("Q266").Value = 3
$ws.Range("R266").Value = "Hortaliza"
